$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.574.93"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "3.790.95"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "708.24"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.81"
$ws.Range("D7").Value = "3.788.91"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  -2.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.17"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "4.430.46"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "3.808.97"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").Value = "70.597.91"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.14"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.35"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "493.57"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.61"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.726"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000144"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.45"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").Value = "3.941.43"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -4.87%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.32"
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.07"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "3.759.12"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.04"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.92"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.27"
$ws.Range("E43").Value = "  -4.79%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "164.55"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.89"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "420.36"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.67"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  -1.54%  "
